$d = $word.ActiveDocument

# Locate the paragraph that contains "Hello" (the first paragraph) using Find
# rather than a hard-coded index, so the script is resilient to minor
# structural differences.
$findRange = $d.Content
$found = $findRange.Find.Execute("Hello")
if (-not $found) {
    throw "Could not find 'Hello' text in the document"
}
$targetParaIndex = $findRange.Paragraphs.Item(1).Index

$p = $d.Paragraphs.Item($targetParaIndex)

# Append " Ritu" (plus a throw-away sentinel character) right before the
# paragraph mark. A trailing sentinel lets us anchor the new bookmark
# exactly at the end of the visible text (immediately after " Ritu") without
# landing on the paragraph-mark boundary itself, which this host mis-handles.
$insertPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$insertPoint.InsertAfter(" RituX")

# Re-fetch the paragraph range (text length changed) and force the newly
# typed " Ritu" text into its own run - distinct from the original "Hello"
# run - by toggling a character format on just that span and back off again.
# (Word normally keeps freshly-typed text in its own run; this host merges
# same-format runs on insert, so nudging formatting forces the run split
# that's present in the target document.)
$p = $d.Paragraphs.Item($targetParaIndex)
$newTextRange = $d.Range($p.Range.End - 7, $p.Range.End - 1)
$newTextRange.Font.Bold = 1
$newTextRange.Font.Bold = 0

# Place the "_GoBack" bookmark right after " Ritu" (still before the sentinel
# "X" and the paragraph mark). Adding a bookmark with a name that already
# exists elsewhere in the document moves it here, removing the old
# occurrence automatically - exactly mirroring the diff, which relocates the
# bookmark from its old paragraph to this one.
$p = $d.Paragraphs.Item($targetParaIndex)
$bookmarkPos = $p.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the sentinel "X" character, leaving "Hello Ritu" followed
# immediately by the bookmark start/end.
$p = $d.Paragraphs.Item($targetParaIndex)
$sentinelRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$sentinelRange.Delete()

Write-Output "Inserted ' Ritu' and relocated the _GoBack bookmark after it."
